# The workbook's only sheet "Síntese" is renamed to "Sintese" (accent
# removed) and the active cell selection moves from E16 to D16.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Sintese"
$ws.Range("D16").Select()
